$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.814.36"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.757.43"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.11"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5071"
$ws.Range("E7").Value = "  +3.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.23"
$ws.Range("E8").Value = "  -1.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2623"
$ws.Range("E9").Value = "  +8.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06190"
$ws.Range("E10").Value = "  +2.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.754.49"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06938"
$ws.Range("E12").Value = "  +4.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.46"
$ws.Range("E13").Value = "  +6.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6022"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.63"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.448"
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.860.19"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006810"
$ws.Range("E20").Value = "  +8.55%  "
$ws.Range("E21").Value = "  +3.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.976.96"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.060"
$ws.Range("E23").Value = "  +4.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.148"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.167"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.79"
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.459"
$ws.Range("E27").Value = "  -2.95%  "
$ws.Range("E28").Value = "  +4.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.808"
$ws.Range("E29").Value = "  -3.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.63"
$ws.Range("E30").Value = "  +3.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08229"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.690"
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("E33").Value = "  +6.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04371"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.656"
$ws.Range("E36").Value = "  +1.72%  "
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.736"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01548"
$ws.Range("E40").Value = "  +4.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.925"
$ws.Range("E41").Value = "  -7.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.31"
$ws.Range("E43").Value = "  +2.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3809"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7436"
$ws.Range("E45").Value = "  -5.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.864"
$ws.Range("E46").Value = "  -5.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05487"
$ws.Range("E47").Value = "  +8.21%  "
$ws.Range("E48").Value = "  +4.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.938"
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.13"
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  +0.29%  "
